# Apply the "streamlined presets mechanism & keeping track of protocol name"
# revision to the Serial_commands worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Reporting section: rename `proto?` command to `p?` and change its
#    reply description from "newline delimited" to "tab delimited".
# ------------------------------------------------------------------
$ws.Range("C8").Value = "Report current protocol info, tab delimited:" + [char]10 + "- protocol name" + [char]10 + "- N_lines"
$ws.Range("B8").Value = "p?"

# ------------------------------------------------------------------
# 2) Control section: the "Additional reply" column for stop/pause/</>/goto
#    now shows the templated reply `{pos}` instead of the old wording.
# ------------------------------------------------------------------
$ws.Range("D14").Value = "{pos}"
$ws.Range("D15").Value = "{pos}"
$ws.Range("D16").Value = "{pos}"
$ws.Range("D17").Value = "{pos}"
$ws.Range("D18").Value = "{pos}"

# ------------------------------------------------------------------
# 3) Debugging section: `b?` now pretty-prints just the line buffer, and a
#    new `proto?` row is inserted right after it that pretty-prints the
#    full protocol program. Everything below shifts down one row.
# ------------------------------------------------------------------
$ws.Rows.Item(27).Insert()
$ws.Range("B26:C26").Copy()
$ws.Range("B27:C27").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("B27").Value = "proto?"
$ws.Range("C27").Value = "Pretty print the full protocol program"
$ws.Range("C26").Value = "Pretty print the current line buffer contents"

# ------------------------------------------------------------------
# 4) Update the print area to match the extra row.
# ------------------------------------------------------------------
$ws.PageSetup.PrintArea = '$B$1:$D$31'

# ------------------------------------------------------------------
# 5) Update the saved view/selection.
# ------------------------------------------------------------------
[void]$ws.Range("F15").Select()
